# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Road ("R") row target depth data ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 466
$wsOff.Range("C3").Value = 295
$wsOff.Range("D3").Value = 78
$wsOff.Range("E3").Value = 35
$wsOff.Range("F3").Value = 10
$wsOff.Range("G3").Value = 4

# --- DEF sheet: update Road ("R") row target depth data ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 460
$wsDef.Range("C3").Value = 346
$wsDef.Range("D3").Value = 89
$wsDef.Range("E3").Value = 39
$wsDef.Range("F3").Value = 2
$wsDef.Range("G3").Value = 2
